$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# RE03 description (B10): add mention of sharing the "endereço eletrônico" (link) of the book
$ws.Range("B10").Value = "Permitir e oferecer o recurso para que quaisquer dos livros visualizados pelo usuário possa ter seu endereço eletrônico compartilhado através do facebook e twitter, independente de ter sido adquirido, e permitir comentário do usuário a respeito da referida obra neste mesmo compartilhamento, a fim de que outras pessoas, ao se interessarem por esta, possam visualizá-la na e-Books Store após devido cadastro (respeitadas e aceitas as políticas de privacidade das respectivas plataformas)."

# RE09 description (B16): "ou aquisições" -> "e de aquisições"
$ws.Range("B16").Value = "Permitir que o comprador visualize online e/ou faça o download da obra adquirida, através de seu histórico de compras e de aquisições gratuitas."

# RE19 description (B26): append "ao registro do evento"
$ws.Range("B26").Value = "Fazer o registro e informar ao cliente o número de protocolo, bem como conferir acesso a este e ao histórico de suas reclamações e respectivas gravações até 90 dias anteriores ao registro do evento."

# Row 10 grows taller to accommodate the longer RE03 text
$ws.Rows.Item(10).RowHeight = 135
